# Insert a new weekly price-record row for "Macroferia Regional de Talca -
# Zanahoria" above the current row 239. This shifts the existing rows
# 239:304 down to 240:305 (preserving all of their data), and the new
# row 239 is populated with this week's reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 239 downward (rows 239:304 -> 240:305) to make room for the
# new record.
$ws.Rows("239:239").Insert()

# Fill in the new row 239 with the new weekly record.
$ws.Range("A239").Value = 5
$ws.Range("B239").Value = "Macroferia Regional de Talca"
$ws.Range("C239").Value = "Maule"
$ws.Range("D239").Value = 44642
$ws.Range("E239").Value = 7
$ws.Range("F239").Value = 100114013
$ws.Range("G239").Value = "Zanahoria"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 500
$ws.Range("K239").Value = 6500
$ws.Range("L239").Value = 6500
$ws.Range("M239").Value = 6500
$ws.Range("N239").Value = "`$/saco 20 kilos"
$ws.Range("O239").Value = "Región de Ñuble"
$ws.Range("P239").Value = 325
$ws.Range("Q239").Value = 20
$ws.Range("R239").Value = "Hortaliza"
